# Update Sheets via scheduled runner
# Applies the numeric corrections captured in the commit diff, sheet by sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H43").Value2 = 2869.2307
$ws.Range("J43").Value2 = 3125
$ws.Range("L43").Value2 = 3125
$ws.Range("N43").Value2 = -3263

$ws.Range("H62").Value2 = 4179.143
$ws.Range("I62").Value2 = 4179.143
$ws.Range("K62").Value2 = 4179.143
$ws.Range("M62").Value2 = -3555.143

$ws.Range("H65").Value2 = 4179.143
$ws.Range("I65").Value2 = 4179.143
$ws.Range("K65").Value2 = 20895.715
$ws.Range("M65").Value2 = -17775.715

$ws.Range("H111").Value2 = 750
$ws.Range("J111").Value2 = 800
$ws.Range("L111").Value2 = 2400
$ws.Range("N111").Value2 = -8534

$ws.Range("H141").Value2 = 4028
$ws.Range("I141").Value2 = 2848.6
$ws.Range("J141").Value2 = 5502.25
$ws.Range("K141").Value2 = 8545.799999999999
$ws.Range("L141").Value2 = 16506.75
$ws.Range("M141").Value2 = -3365.799999999999
$ws.Range("N141").Value2 = -26866.75

# ---------------------------------------------------------------
# ARM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H44").Value2 = 57921.145
$ws.Range("J44").Value2 = 57921.145
$ws.Range("L44").Value2 = 57921.145
$ws.Range("N44").Value2 = -58897.145

$ws.Range("H45").Value2 = 797.41
$ws.Range("I45").Value2 = 797.5859
$ws.Range("J45").Value2 = 780
$ws.Range("K45").Value2 = 797.5859
$ws.Range("L45").Value2 = 780
$ws.Range("M45").Value2 = -420.5859
$ws.Range("N45").Value2 = -1534

$ws.Range("H61").Value2 = 2475.6785
$ws.Range("I61").Value2 = 1942.9333
$ws.Range("J61").Value2 = 3090.3845
$ws.Range("K61").Value2 = 1942.9333
$ws.Range("L61").Value2 = 3090.3845
$ws.Range("M61").Value2 = -1730.9333
$ws.Range("N61").Value2 = -3514.3845

$ws.Range("H74").Value2 = 1389.4897
$ws.Range("I74").Value2 = 1153.7391
$ws.Range("J74").Value2 = 5004.3335
$ws.Range("K74").Value2 = 1153.7391
$ws.Range("L74").Value2 = 5004.3335
$ws.Range("M74").Value2 = -279.7391
$ws.Range("N74").Value2 = -6752.3335

$ws.Range("H77").Value2 = 1389.4897
$ws.Range("I77").Value2 = 1153.7391
$ws.Range("J77").Value2 = 5004.3335
$ws.Range("K77").Value2 = 5768.6955
$ws.Range("L77").Value2 = 25021.6675
$ws.Range("M77").Value2 = -1400.6955
$ws.Range("N77").Value2 = -33757.6675

$ws.Range("H122").Value2 = 2008.6154
$ws.Range("I122").Value2 = 2008.6154
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 6025.8462
$ws.Range("L122").Value2 = 0
$ws.Range("M122").Value2 = -3575.8462
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value2 = 13516467
$ws.Range("I132").Value2 = 22729036
$ws.Range("J132").Value2 = 4698.067
$ws.Range("K132").Value2 = 68187108
$ws.Range("L132").Value2 = 14094.201
$ws.Range("M132").Value2 = -68184578
$ws.Range("N132").Value2 = -19154.201

$ws.Range("H136").Value2 = 2475.6785
$ws.Range("I136").Value2 = 1942.9333
$ws.Range("J136").Value2 = 3090.3845
$ws.Range("K136").Value2 = 5828.7999
$ws.Range("L136").Value2 = 9271.1535
$ws.Range("M136").Value2 = -3278.7999
$ws.Range("N136").Value2 = -14371.1535

# ---------------------------------------------------------------
# BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value2 = 4568.943
$ws.Range("I20").Value2 = 1133
$ws.Range("J20").Value2 = 7145.9
$ws.Range("K20").Value2 = 1133
$ws.Range("L20").Value2 = 7145.9
$ws.Range("M20").Value2 = -886
$ws.Range("N20").Value2 = -7639.9

$ws.Range("H21").Value2 = 19886.4
$ws.Range("J21").Value2 = 19886.4
$ws.Range("L21").Value2 = 19886.4
$ws.Range("N21").Value2 = -20358.4

$ws.Range("H86").Value2 = 5899.75
$ws.Range("I86").Value2 = 4533.3335
$ws.Range("J86").Value2 = 9999
$ws.Range("K86").Value2 = 4533.3335
$ws.Range("L86").Value2 = 9999
$ws.Range("M86").Value2 = -3410.3335
$ws.Range("N86").Value2 = -12245

$ws.Range("H89").Value2 = 5899.75
$ws.Range("I89").Value2 = 4533.3335
$ws.Range("J89").Value2 = 9999
$ws.Range("K89").Value2 = 22666.6675
$ws.Range("L89").Value2 = 49995
$ws.Range("M89").Value2 = -17050.6675
$ws.Range("N89").Value2 = -61227

$ws.Range("H99").Value2 = 2278.7917
$ws.Range("I99").Value2 = 2193.5293
$ws.Range("J99").Value2 = 2485.8572
$ws.Range("K99").Value2 = 2193.5293
$ws.Range("L99").Value2 = 2485.8572
$ws.Range("M99").Value2 = -695.5293000000001
$ws.Range("N99").Value2 = -5481.8572

$ws.Range("H107").Value2 = 3768.125
$ws.Range("I107").Value2 = 2876
$ws.Range("J107").Value2 = 10013
$ws.Range("K107").Value2 = 2876
$ws.Range("L107").Value2 = 10013
$ws.Range("M107").Value2 = -956
$ws.Range("N107").Value2 = -13853

$ws.Range("H115").Value2 = 19999.889
$ws.Range("J115").Value2 = 19999.889
$ws.Range("L115").Value2 = 19999.889
$ws.Range("N115").Value2 = -23133.889

$ws.Range("H134").Value2 = 1638.1515
$ws.Range("I134").Value2 = 1139.3704
$ws.Range("J134").Value2 = 3882.6667
$ws.Range("K134").Value2 = 3418.1112
$ws.Range("L134").Value2 = 11648.0001
$ws.Range("M134").Value2 = -883.1112000000003
$ws.Range("N134").Value2 = -16718.0001

# ---------------------------------------------------------------
# CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H48").Value2 = 35000
$ws.Range("I48").Value2 = 0
$ws.Range("K48").Value2 = 0
$ws.Range("M48").ClearContents()

$ws.Range("H134").Value2 = 935630.4399999999
$ws.Range("I134").Value2 = 1906
$ws.Range("J134").Value2 = 2803079.5
$ws.Range("K134").Value2 = 5718
$ws.Range("L134").Value2 = 8409238.5
$ws.Range("M134").Value2 = -3183
$ws.Range("N134").Value2 = -8414308.5

# ---------------------------------------------------------------
# GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H132").Value2 = 2740.9546
$ws.Range("I132").Value2 = 2181.875
$ws.Range("J132").Value2 = 4231.8335
$ws.Range("K132").Value2 = 6545.625
$ws.Range("L132").Value2 = 12695.5005
$ws.Range("M132").Value2 = -4015.625
$ws.Range("N132").Value2 = -17755.5005

$ws.Range("H140").Value2 = 38168.4
$ws.Range("J140").Value2 = 38168.4
$ws.Range("L140").Value2 = 38168.4
$ws.Range("N140").Value2 = -48528.4

$ws.Range("H141").Value2 = 70104.75
$ws.Range("J141").Value2 = 70104.75
$ws.Range("L141").Value2 = 70104.75
$ws.Range("N141").Value2 = -80464.75

# ---------------------------------------------------------------
# LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H45").Value2 = 18890
$ws.Range("I45").Value2 = 0
$ws.Range("J45").Value2 = 18890
$ws.Range("K45").Value2 = 0
$ws.Range("L45").Value2 = 18890
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value2 = -19704

$ws.Range("H53").Value2 = 29975
$ws.Range("J53").Value2 = 29975
$ws.Range("L53").Value2 = 29975
$ws.Range("N53").Value2 = -31011

$ws.Range("H82").Value2 = 10421102
$ws.Range("I82").Value2 = 3994
$ws.Range("J82").Value2 = 16671367
$ws.Range("K82").Value2 = 3994
$ws.Range("L82").Value2 = 16671367
$ws.Range("M82").Value2 = -3633
$ws.Range("N82").Value2 = -16672089

$ws.Range("H85").Value2 = 10421102
$ws.Range("I85").Value2 = 3994
$ws.Range("J85").Value2 = 16671367
$ws.Range("K85").Value2 = 3994
$ws.Range("L85").Value2 = 16671367
$ws.Range("M85").Value2 = -2746
$ws.Range("N85").Value2 = -16673863

$ws.Range("H93").Value2 = 1093.375
$ws.Range("I93").Value2 = 665
$ws.Range("J93").Value2 = 1192.2307
$ws.Range("K93").Value2 = 665
$ws.Range("L93").Value2 = 1192.2307
$ws.Range("M93").Value2 = 583
$ws.Range("N93").Value2 = -3688.2307

$ws.Range("H137").Value2 = 49144.547
$ws.Range("J137").Value2 = 49144.547
$ws.Range("L137").Value2 = 49144.547
$ws.Range("N137").Value2 = -59344.547

# ---------------------------------------------------------------
# WVR
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H53").Value2 = 19495
$ws.Range("J53").Value2 = 19495
$ws.Range("L53").Value2 = 19495
$ws.Range("N53").Value2 = -20709

$ws.Range("H96").Value2 = 2100.4285
$ws.Range("I96").Value2 = 2183.8333
$ws.Range("J96").Value2 = 1600
$ws.Range("K96").Value2 = 2183.8333
$ws.Range("L96").Value2 = 1600
$ws.Range("M96").Value2 = -810.8332999999998
$ws.Range("N96").Value2 = -4346

$ws.Range("H107").Value2 = 6250660
$ws.Range("I107").Value2 = 635.36365
$ws.Range("J107").Value2 = 20000714
$ws.Range("K107").Value2 = 1906.09095
$ws.Range("L107").Value2 = 60002142
$ws.Range("M107").Value2 = 13.90904999999998
$ws.Range("N107").Value2 = -60005982

$ws.Range("H126").Value2 = 4203438
$ws.Range("I126").Value2 = 4203438
$ws.Range("K126").Value2 = 12610314
$ws.Range("M126").Value2 = -12607844

$ws.Range("H136").Value2 = 213696.28
$ws.Range("I136").Value2 = 257131.33
$ws.Range("J136").Value2 = 1950.375
$ws.Range("K136").Value2 = 771393.99
$ws.Range("L136").Value2 = 5851.125
$ws.Range("M136").Value2 = -768843.99
$ws.Range("N136").Value2 = -10951.125

$ws.Range("H140").Value2 = 33981.582
$ws.Range("J140").Value2 = 33981.582
$ws.Range("L140").Value2 = 33981.582
$ws.Range("N140").Value2 = -44341.582

$ws.Range("H141").Value2 = 38913.715
$ws.Range("J141").Value2 = 38913.715
$ws.Range("L141").Value2 = 38913.715
$ws.Range("N141").Value2 = -49273.715
